$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-06-18 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-06-19 Thursday", 2)

$d.Content.Find.Execute("140×6=", $true, $false, $false, $false, $false, $true, 1, $false, "720×7=", 2)
$d.Content.Find.Execute("539×4=", $true, $false, $false, $false, $false, $true, 1, $false, "617×8=", 2)
$d.Content.Find.Execute("201×3=", $true, $false, $false, $false, $false, $true, 1, $false, "920×2=", 2)
$d.Content.Find.Execute("118×8=", $true, $false, $false, $false, $false, $true, 1, $false, "821×5=", 2)
$d.Content.Find.Execute("478×3=", $true, $false, $false, $false, $false, $true, 1, $false, "233×4=", 2)

$d.Content.Find.Execute("444×7=", $true, $false, $false, $false, $false, $true, 1, $false, "396×9=", 2)
$d.Content.Find.Execute("283×9=", $true, $false, $false, $false, $false, $true, 1, $false, "252×2=", 2)
$d.Content.Find.Execute("959×8=", $true, $false, $false, $false, $false, $true, 1, $false, "305×2=", 2)
$d.Content.Find.Execute("565×6=", $true, $false, $false, $false, $false, $true, 1, $false, "709×4=", 2)
$d.Content.Find.Execute("861×3=", $true, $false, $false, $false, $false, $true, 1, $false, "186×9=", 2)

$d.Content.Find.Execute("299×5=", $true, $false, $false, $false, $false, $true, 1, $false, "214×7=", 2)
$d.Content.Find.Execute("421×3=", $true, $false, $false, $false, $false, $true, 1, $false, "921×8=", 2)
$d.Content.Find.Execute("414×3=", $true, $false, $false, $false, $false, $true, 1, $false, "115×4=", 2)
$d.Content.Find.Execute("288×3=", $true, $false, $false, $false, $false, $true, 1, $false, "624×9=", 2)
$d.Content.Find.Execute("988×2=", $true, $false, $false, $false, $false, $true, 1, $false, "578×7=", 2)

$d.Content.Find.Execute("577×6=", $true, $false, $false, $false, $false, $true, 1, $false, "738×3=", 2)
$d.Content.Find.Execute("723×2=", $true, $false, $false, $false, $false, $true, 1, $false, "328×5=", 2)
$d.Content.Find.Execute("266×2=", $true, $false, $false, $false, $false, $true, 1, $false, "723×5=", 2)
$d.Content.Find.Execute("220×2=", $true, $false, $false, $false, $false, $true, 1, $false, "620×6=", 2)
$d.Content.Find.Execute("828×6=", $true, $false, $false, $false, $false, $true, 1, $false, "988×8=", 2)

$d.Content.Find.Execute("981×9=", $true, $false, $false, $false, $false, $true, 1, $false, "650×3=", 2)
$d.Content.Find.Execute("610×4=", $true, $false, $false, $false, $false, $true, 1, $false, "949×8=", 2)
$d.Content.Find.Execute("411×8=", $true, $false, $false, $false, $false, $true, 1, $false, "308×9=", 2)
$d.Content.Find.Execute("530×9=", $true, $false, $false, $false, $false, $true, 1, $false, "555×8=", 2)
$d.Content.Find.Execute("109×3=", $true, $false, $false, $false, $false, $true, 1, $false, "708×8=", 2)
